$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 200, shifting existing rows 200..258 down to 201..259
$ws.Rows.Item(200).Insert()

# Populate the newly inserted row 200 with data
$ws.Cells.Item(200, 1).Value = 10
$ws.Cells.Item(200, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(200, 3).Value = "La Araucanía"
$ws.Cells.Item(200, 4).Value = 44627
$ws.Cells.Item(200, 5).Value = 9
$ws.Cells.Item(200, 6).Value = 100112017
$ws.Cells.Item(200, 7).Value = "Apio"
$ws.Cells.Item(200, 8).Value = "Americana (o)"
$ws.Cells.Item(200, 9).Value = "Primera"
$ws.Cells.Item(200, 10).Value = 65
$ws.Cells.Item(200, 11).Value = 8000
$ws.Cells.Item(200, 12).Value = 8000
$ws.Cells.Item(200, 13).Value = 8000
$ws.Cells.Item(200, 14).Value = "`$/docena de matas"
$ws.Cells.Item(200, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(200, 16).Value = 1333
$ws.Cells.Item(200, 17).Value = 6
$ws.Cells.Item(200, 18).Value = "Hortaliza"
